# no-op for now, just to see default re-save behavior
$p = $ppt.ActivePresentation
